# Update games of 2019-06-01
$wb = $excel.ActiveWorkbook

$wsResultats = $wb.Worksheets.Item("Résultats")
$wsPool = $wb.Worksheets.Item("Pool")

# Enter game 2 and game 3 scores for the BOSTON vs ST-LOUIS series
$wsResultats.Range("AB15").Value = 2
$wsResultats.Range("AC15").Value = 7
$wsResultats.Range("AB16").Value = 3
$wsResultats.Range("AC16").Value = 2

# Update the selection/scroll position on the Résultats sheet
$wsResultats.Activate()
$wsResultats.Range("AB16").Select()

# Update the selection/scroll position on the Pool sheet (this sheet stays active)
$wsPool.Activate()
$wsPool.Application.ActiveWindow.ScrollRow = 7
$wsPool.Range("B35").Select()
